# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The worker-detail table (rows 16-38, cols B:G) on "Hoja1" is rebuilt:
# the old account-statement periods for SILFREDO MARTINEZ CABARCAS
# (2207..2404, ascending) are replaced with the same periods listed in
# DESCENDING order, and the single ADOL ANTONIO TORRES TRESPALACIOS
# row (period 2111) moves from the top of the block (row 16) to the
# bottom (row 38) carrying its own F/G values with it, while row 16
# picks up the F/G values that used to belong to the old "2404" entry.
#
# Net effect per row (B/C/D/E/F/G), row-by-row:
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 16; B = "CC"; C = "1126118302"; D = "ADOL ANTONIO TORRES TRESPALACIOS";  E = "2404"; F = 46000; G = 1500000 },
    @{ Row = 17; B = "CC"; C = "73136692";   D = "SILFREDO MARTINEZ CABARCAS";        E = "2403"; F = 60000; G = 1500000 },
    @{ Row = 18; B = "CC"; C = "73136692";   D = "SILFREDO MARTINEZ CABARCAS";        E = "2402"; F = 60000; G = 1500000 },
    @{ Row = 19; B = "CC"; C = "73136692";   D = "SILFREDO MARTINEZ CABARCAS";        E = "2401"; F = 60000; G = 1500000 },
    @{ Row = 20; B = "CC"; C = "73136692";   D = "SILFREDO MARTINEZ CABARCAS";        E = "2312"; F = 60000; G = 1500000 },
    @{ Row = 21; B = "CC"; C = "73136692";   D = "SILFREDO MARTINEZ CABARCAS";        E = "2311"; F = 60000; G = 1500000 },
    @{ Row = 22; B = "CC"; C = "73136692";   D = "SILFREDO MARTINEZ CABARCAS";        E = "2310"; F = 60000; G = 1500000 },
    @{ Row = 23; B = "CC"; C = "73136692";   D = "SILFREDO MARTINEZ CABARCAS";        E = "2309"; F = 60000; G = 1500000 },
    @{ Row = 24; B = "CC"; C = "73136692";   D = "SILFREDO MARTINEZ CABARCAS";        E = "2308"; F = 60000; G = 1500000 },
    @{ Row = 25; B = "CC"; C = "73136692";   D = "SILFREDO MARTINEZ CABARCAS";        E = "2307"; F = 60000; G = 1500000 },
    @{ Row = 26; B = "CC"; C = "73136692";   D = "SILFREDO MARTINEZ CABARCAS";        E = "2306"; F = 60000; G = 1500000 },
    @{ Row = 27; B = "CC"; C = "73136692";   D = "SILFREDO MARTINEZ CABARCAS";        E = "2305"; F = 60000; G = 1500000 },
    @{ Row = 28; B = "CC"; C = "73136692";   D = "SILFREDO MARTINEZ CABARCAS";        E = "2304"; F = 60000; G = 1500000 },
    @{ Row = 29; B = "CC"; C = "73136692";   D = "SILFREDO MARTINEZ CABARCAS";        E = "2303"; F = 60000; G = 1500000 },
    @{ Row = 30; B = "CC"; C = "73136692";   D = "SILFREDO MARTINEZ CABARCAS";        E = "2302"; F = 60000; G = 1500000 },
    @{ Row = 31; B = "CC"; C = "73136692";   D = "SILFREDO MARTINEZ CABARCAS";        E = "2301"; F = 60000; G = 1500000 },
    @{ Row = 32; B = "CC"; C = "73136692";   D = "SILFREDO MARTINEZ CABARCAS";        E = "2212"; F = 60000; G = 1500000 },
    @{ Row = 33; B = "CC"; C = "73136692";   D = "SILFREDO MARTINEZ CABARCAS";        E = "2211"; F = 60000; G = 1500000 },
    @{ Row = 34; B = "CC"; C = "73136692";   D = "SILFREDO MARTINEZ CABARCAS";        E = "2210"; F = 60000; G = 1500000 },
    @{ Row = 35; B = "CC"; C = "73136692";   D = "SILFREDO MARTINEZ CABARCAS";        E = "2209"; F = 60000; G = 1500000 },
    @{ Row = 36; B = "CC"; C = "73136692";   D = "SILFREDO MARTINEZ CABARCAS";        E = "2208"; F = 60000; G = 1500000 },
    @{ Row = 37; B = "CC"; C = "73136692";   D = "SILFREDO MARTINEZ CABARCAS";        E = "2207"; F = 24000; G = 1500000 },
    @{ Row = 38; B = "CC"; C = "1126118302"; D = "ADOL ANTONIO TORRES TRESPALACIOS";  E = "2111"; F = 30666; G = 1000000 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("B$n").Value = $r.B
    $ws.Range("C$n").Value = $r.C
    $ws.Range("D$n").Value = $r.D
    $ws.Range("E$n").Value = $r.E
    $ws.Range("F$n").Value = $r.F
    $ws.Range("G$n").Value = $r.G
}
